$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 8.800000000000001
$ws.Range("I3").Value = 12
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.81
$ws.Range("O3").Value = 1.34
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.26
$ws.Range("S3").Value = 3.05
$ws.Range("T3").Value = 1.87
$ws.Range("U3").Value = 1.44
$ws.Range("V3").Value = 1.09
$ws.Range("W3").Value = 3
$ws.Range("X3:AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 1.93
$ws.Range("K4").Value = 5.4
$ws.Range("P4").Value = 2.68
$ws.Range("Q4").Value = 1.41

# Row 5
$ws.Range("Q5").Value = 1.29

# Row 6
$ws.Range("J6").Value = 3.9

# Row 7
$ws.Range("J7").Value = 4.4

# Row 8
$ws.Range("F8").Value = 2.74
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 2.54
$ws.Range("J8").Value = 3.35
